$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.138873100280762
$ws.Range("B1").Value = 6.007148265838623
$ws.Range("C1").Value = 2.4826340675354
$ws.Range("D1").Value = 1.128854632377625
$ws.Range("E1").Value = 0.8033167123794556
